$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted above the current first data row (288),
# pushing all existing records (288-304) down by one (to 289-305).
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new weekly price record.
$ws.Cells.Item(288, 1).Value = 4
$ws.Cells.Item(288, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(288, 3).Value = "Los Lagos"
$ws.Cells.Item(288, 4).Value = 44783
$ws.Cells.Item(288, 5).Value = 10
$ws.Cells.Item(288, 6).Value = 100112043
$ws.Cells.Item(288, 7).Value = "Pepino ensalada"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 70
$ws.Cells.Item(288, 11).Value = 27000
$ws.Cells.Item(288, 12).Value = 27000
$ws.Cells.Item(288, 13).Value = 27000
$ws.Cells.Item(288, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(288, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(288, 16).Value = 450
$ws.Cells.Item(288, 17).Value = 60
$ws.Cells.Item(288, 18).Value = "Hortaliza"
